$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values such as "206.14" or "1.698.92" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.867.75"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.563.83"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "206.14"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "21.78"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "0.0585"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.785.35"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.564.26"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "26.864.51"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "61.25"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "215.27"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "9.19"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "153.60"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "6.74"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "1.405.52"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "0.920"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "0.530"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").Value = "0.811"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").Value = "2.18"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "63.42"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "1.698.92"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "86.35"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").Value = "0.0₇0976"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  +0.71%  "
